$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Feria Lagunitas de Puerto Montt -
# Acelga". It belongs at row 220, pushing the existing rows 220-271 down to
# 221-272 (dimension grows from A1:R271 to A1:R272).
$ws.Rows.Item(220).Insert()

# Populate the newly inserted row 220 with the new record.
$ws.Cells.Item(220, 1).Value = 4
$ws.Cells.Item(220, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(220, 3).Value = "Los Lagos"
$ws.Cells.Item(220, 4).Value = 44964
$ws.Cells.Item(220, 5).Value = 10
$ws.Cells.Item(220, 6).Value = 100112009
$ws.Cells.Item(220, 7).Value = "Acelga"
$ws.Cells.Item(220, 8).Value = "Sin especificar"
$ws.Cells.Item(220, 9).Value = "Primera"
$ws.Cells.Item(220, 10).Value = 80
$ws.Cells.Item(220, 11).Value = 10000
$ws.Cells.Item(220, 12).Value = 10000
$ws.Cells.Item(220, 13).Value = 10000
$ws.Cells.Item(220, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(220, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(220, 16).Value = 833
$ws.Cells.Item(220, 17).Value = 12
$ws.Cells.Item(220, 18).Value = "Hortaliza"
